$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.274.85'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '2.965.18'
$ws.Range('D5').Value = '587.14'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').Value = '141.68'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -2.69%  '
$ws.Range('D9').Value = '2.963.08'
$ws.Range('E9').Value = '  -2.89%  '
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -6.53%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '0.453'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('D13').Value = '0.0000224'
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').Value = '33.85'
$ws.Range('E14').Value = '  -5.48%  '
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '3.457.40'
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').Value = '61.277.63'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('D19').Value = '2.964.94'
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('D20').Value = '446.25'
$ws.Range('E20').Value = '  -6.67%  '
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('D22').Value = '0.681'
$ws.Range('E22').Value = '  -3.21%  '
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('D24').Value = '81.10'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('E25').Value = '  -4.14%  '
$ws.Range('D26').Value = '2.16'
$ws.Range('E26').Value = '  -9.25%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  -5.91%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').Value = '6.82'
$ws.Range('E31').Value = '  -5.95%  '
$ws.Range('E32').Value = '  -6.44%  '
$ws.Range('D33').Value = '26.87'
$ws.Range('E33').Value = '  -3.12%  '
$ws.Range('E34').Value = '  -4.00%  '
$ws.Range('E35').Value = '  -4.59%  '
$ws.Range('D36').Value = '0.0₃0776'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('D37').Value = '5.70'
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('D38').Value = '50.14'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  -5.95%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('E41').Value = '  +4.54%  '
$ws.Range('D42').Value = '2.75'
$ws.Range('E42').Value = '  -10.14%  '
$ws.Range('D43').Value = '387.77'
$ws.Range('E43').Value = '  -8.80%  '
$ws.Range('D44').Value = '0.0351'
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('D45').Value = '0.263'
$ws.Range('E45').Value = '  -7.59%  '
$ws.Range('D46').Value = '2.685.22'
$ws.Range('E46').Value = '  -5.18%  '
$ws.Range('D47').Value = '37.04'
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').Value = '130.82'
$ws.Range('E48').Value = '  +2.77%  '
$ws.Range('D50').Value = '0.108'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('E51').Value = '  -1.31%  '
